$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.859.40"
$ws.Range("E2").Value = "  -6.33%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.442.31"
$ws.Range("E3").Value = "  -8.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "538.67"
$ws.Range("E5").Value = "  -2.71%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'146.60"
$ws.Range("E6").Value = "  -7.22%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.11%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -3.10%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.454.89"
$ws.Range("E9").Value = "  -8.55%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0988"
$ws.Range("E10").Value = "  -6.74%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -2.18%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  -1.23%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "'0.350"
$ws.Range("E13").Value = "  -5.05%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.880.60"
$ws.Range("E14").Value = "  -8.68%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "23.88"
$ws.Range("E15").Value = "  -9.96%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "58.750.42"
$ws.Range("E16").Value = "  -6.37%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -6.23%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.505.22"
$ws.Range("E18").Value = "  -6.61%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -6.52%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "4.35"
$ws.Range("E20").Value = "  -5.63%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "322.84"
$ws.Range("E21").Value = "  -6.40%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -3.37%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.68"
$ws.Range("E23").Value = "  -9.69%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "60.59"
$ws.Range("E24").Value = "  -3.98%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  -11.21%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -5.24%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").Value = "0.976"
$ws.Range("E27").Value = "  -2.33%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  -6.49%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -6.13%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0766"
$ws.Range("E30").Value = "  -10.50%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -12.56%  "

# Row 32 - Aptos
$ws.Range("D32").Value = "6.63"
$ws.Range("E32").Value = "  -8.43%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.06%  "

# Row 34 - Monero
$ws.Range("D34").Value = "'156.20"
$ws.Range("E34").Value = "  -4.87%  "

# Row 35 - swapped to ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.37"
$ws.Range("E35").Value = "  -7.50%  "

# Row 36 - swapped to EthereumClassic
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "18.38"
$ws.Range("E36").Value = "  -5.62%  "

# Row 37 - NEARProtocol
$ws.Range("D37").Value = "4.44"
$ws.Range("E37").Value = "  -9.88%  "

# Row 38 - Stacks
$ws.Range("D38").Value = "1.69"
$ws.Range("E38").Value = "  -5.04%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  -6.50%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "312.25"
$ws.Range("E40").Value = "  -10.33%  "

# Row 41 - OKB
$ws.Range("D41").Value = "36.16"
$ws.Range("E41").Value = "  -5.65%  "

# Row 42 - SuiNetwork
$ws.Range("D42").Value = "0.831"
$ws.Range("E42").Value = "  -12.16%  "

# Row 43 - Filecoin
$ws.Range("D43").Value = "3.69"
$ws.Range("E43").Value = "  -7.35%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.22%  "

# Row 45 - WhiteBITCoin
$ws.Range("D45").Value = "10.74"
$ws.Range("E45").Value = "  -2.36%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  -2.96%  "

# Row 48 - Hedera
$ws.Range("D48").Value = "0.0523"
$ws.Range("E48").Value = "  -6.28%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  -5.23%  "

# Row 50 - Aave
$ws.Range("D50").Value = "121.75"
$ws.Range("E50").Value = "  -5.60%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "18.77"
$ws.Range("E51").Value = "  -10.21%  "
